# Dodano polska wersje w quizie
# Translate the quiz questions/answers on Arkusz1 from English to Polish.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- "Which country has the highest poverty level?" block (rows 2-5) ---
$ws.Range("D2").Value = "Który kraj ma najwyzszy poziom wskaźnika biedy?"
$ws.Range("D3").Value = "Który kraj ma najwyzszy poziom wskaźnika biedy?"
$ws.Range("E3").Value = "Portugalia"
$ws.Range("D4").Value = "Który kraj ma najwyzszy poziom wskaźnika biedy?"
$ws.Range("D5").Value = "Który kraj ma najwyzszy poziom wskaźnika biedy?"
$ws.Range("E5").Value = "Etiopia"

# --- "Which sentance is true?" block (rows 6-9) ---
$ws.Range("D6").Value = "Które zdanie jest prawdziwe?"
$ws.Range("D7").Value = "Które zdanie jest prawdziwe?"
$ws.Range("D8").Value = "Które zdanie jest prawdziwe?"
$ws.Range("D9").Value = "Które zdanie jest prawdziwe?"

# --- "Which organisation uses this logo?" block (rows 10-13) ---
$ws.Range("D10").Value = "Jakiej organizacji jest to logo? "
$ws.Range("D11").Value = "Jakiej organizacji jest to logo? "
$ws.Range("D12").Value = "Jakiej organizacji jest to logo? "
$ws.Range("D13").Value = "Jakiej organizacji jest to logo? "

# --- "In which country women held the biggest number of seats?" block (rows 14-17) ---
$ws.Range("D14").Value = "W którym kraju kobiety zajmują najwięcej miejsc w parlamencie?"
$ws.Range("E14").Value = "Kanada"
$ws.Range("D15").Value = "W którym kraju kobiety zajmują najwięcej miejsc w parlamencie?"
$ws.Range("E15").Value = "Francja"
$ws.Range("D16").Value = "W którym kraju kobiety zajmują najwięcej miejsc w parlamencie?"
$ws.Range("E16").Value = "Dania"
$ws.Range("D17").Value = "W którym kraju kobiety zajmują najwięcej miejsc w parlamencie?"
$ws.Range("E17").Value = "Czechy "

# --- "How many SDG there are?" block (rows 18-23) ---
$ws.Range("D18").Value = "Ile istnieje celów zwrównoważonego rozwoju? "
$ws.Range("D19").Value = "Ile istnieje celów zwrównoważonego rozwoju? "
$ws.Range("D20").Value = "Ile istnieje celów zwrównoważonego rozwoju? "
$ws.Range("D21").Value = "Ile istnieje celów zwrównoważonego rozwoju? "
$ws.Range("D22").Value = "Ile istnieje celów zwrównoważonego rozwoju? "
$ws.Range("D23").Value = "Ile istnieje celów zwrównoważonego rozwoju? "

# --- "Which SDG Logo you saw?" block (rows 24-31) ---
$ws.Range("D24").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E24").Value = "Czysta i dostępna energia"
$ws.Range("D25").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E25").Value = "Życie na lądzie"
$ws.Range("D26").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E26").Value = "Zero głodu"
$ws.Range("D27").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E27").Value = "Koniec z ubstwem "
$ws.Range("D28").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E28").Value = "Działania w dziedzinie klimatu "
$ws.Range("D29").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E29").Value = "Odpowiedzialna konsumpcja i produkcja"
$ws.Range("D30").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E30").Value = "Życie pod wodą"
$ws.Range("D31").Value = "Jakiego celu logo zobaczyłeś "
$ws.Range("E31").Value = "Dobra jakość edukacji"

# Move the active selection the way the author left it in the saved file.
$ws.Range("E26").Select()
